$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the 2020 column (N) into the new 2021 column (O)
$ws.Range("N4:N5").Copy()
$ws.Range("O4:O5").PasteSpecial(-4122)

# Fill in the new 2021 data
$ws.Range("O4").Value = 2021
$ws.Range("O5").Value = 1.5020015556876996

# Move / update the current selection to match the saved workbook state
$ws.Range("Q5").Select()
